$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a duplicate of the first 48 "NN / 48 miast" solution rows
# (rows 2-49) to the end of the sheet (rows 146-193), matching the
# target workbook's extended dataset.
$source = $ws.Range("A2:D49")
$target = $ws.Range("A146:D193")
$target.Value = $source.Value()

# Writing the multi-line route text triggers Excel's implicit row
# autofit (customHeight). Re-running AutoFit clears the stale/custom
# height flag so the new rows end up with the same "no explicit
# height" row XML as every other data row in the sheet.
$target.Rows.AutoFit()
